$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the existing "_GoBack" bookmark. In the original document it sits
#    in the (empty, bold) paragraph right after "Defeitos gerais:" - it will
#    be re-created at the end of the title paragraph below.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2. Append a new run "– P4" to the end of the title paragraph (paragraph 1).
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titleRange = $titlePara.Range.Duplicate()
$insertPoint = $titleRange.End - 1
$endRange = $d.Range($insertPoint, $insertPoint)
$endRange.InsertAfter("– P4")

# ---------------------------------------------------------------------------
# 3. Re-insert the "_GoBack" bookmark right after the new run, i.e. at the
#    (new) end of the title paragraph.
#
#    Note: placing a bookmark exactly at "paragraph.Range.End - 1" behaves
#    unreliably, so a temporary marker character is appended first, the
#    bookmark is placed just before it (a safe, non-edge position), and the
#    marker is then removed again.
# ---------------------------------------------------------------------------
$titlePara2 = $d.Paragraphs(1)
$titleRange2 = $titlePara2.Range.Duplicate()
$tempPoint = $titleRange2.End - 1
$tempRange = $d.Range($tempPoint, $tempPoint)
$tempRange.InsertAfter("Z")

$titlePara3 = $d.Paragraphs(1)
$titleRange3 = $titlePara3.Range.Duplicate()
$bookmarkPoint = $titleRange3.End - 2
$bookmarkRange = $d.Range($bookmarkPoint, $bookmarkPoint)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$titlePara4 = $d.Paragraphs(1)
$titleRange4 = $titlePara4.Range.Duplicate()
$markerPoint = $titleRange4.End - 2
$markerRange = $d.Range($markerPoint, $markerPoint + 1)
$markerRange.Delete()

# ---------------------------------------------------------------------------
# 4. Add a "1" run into the first (defect-number) cell of the first data row
#    of the table, i.e. row 2, column 1 (row 1 is the "Defeito"/"Descrição"
#    header row).
# ---------------------------------------------------------------------------
$table = $d.Tables(1)
$cell = $table.Cell(2, 1)
$cell.Range.InsertAfter("1")
